# Auto-generated Excel COM-interop script
# Applies numeric corrections to Leve profit-calculation columns (H-N)
# across multiple worksheets, per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 54997.223
$ws.Range("I21").Value = 51871.875
$ws.Range("J21").Value = 80000
$ws.Range("K21").Value = 51871.875
$ws.Range("L21").Value = 80000
$ws.Range("M21").Value = -51403.875
$ws.Range("N21").Value = -80936
$ws.Range("H23").Value = 54997.223
$ws.Range("I23").Value = 51871.875
$ws.Range("J23").Value = 80000
$ws.Range("K23").Value = 51871.875
$ws.Range("L23").Value = 80000
$ws.Range("M23").Value = -51637.875
$ws.Range("N23").Value = -80468
$ws.Range("H38").Value = 1129.4546
$ws.Range("I38").Value = 1162.4
$ws.Range("J38").Value = 800
$ws.Range("K38").Value = 3487.2
$ws.Range("L38").Value = 2400
$ws.Range("M38").Value = -3115.2
$ws.Range("N38").Value = -3144
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 2000
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 2000
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -1064
$ws.Range("N74").Value = -5872
$ws.Range("H76").Value = 6201.636
$ws.Range("I76").Value = 4643.6
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 4643.6
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -4328.6
$ws.Range("N76").Value = -8130
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 2000
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 10000
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -5320
$ws.Range("N77").Value = -29360
$ws.Range("H79").Value = 6201.636
$ws.Range("I79").Value = 4643.6
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 4643.6
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -3551.6
$ws.Range("N79").Value = -9684
$ws.Range("H138").Value = 2087.8572
$ws.Range("I138").Value = 692.86664
$ws.Range("J138").Value = 2773.918
$ws.Range("K138").Value = 2078.59992
$ws.Range("L138").Value = 8321.754000000001
$ws.Range("M138").Value = 3061.40008
$ws.Range("N138").Value = -18601.754

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 629.125
$ws.Range("I2").Value = 510
$ws.Range("J2").Value = 1224.75
$ws.Range("K2").Value = 510
$ws.Range("L2").Value = 1224.75
$ws.Range("M2").Value = -397
$ws.Range("N2").Value = -1450.75
$ws.Range("H45").Value = 22729618
$ws.Range("I45").Value = 35715630
$ws.Range("J45").Value = 4090.5
$ws.Range("K45").Value = 35715630
$ws.Range("L45").Value = 4090.5
$ws.Range("M45").Value = -35715253
$ws.Range("N45").Value = -4844.5
$ws.Range("H61").Value = 14710648
$ws.Range("I61").Value = 11115532
$ws.Range("J61").Value = 41674024
$ws.Range("K61").Value = 11115532
$ws.Range("L61").Value = 41674024
$ws.Range("M61").Value = -11115320
$ws.Range("N61").Value = -41674448
$ws.Range("H74").Value = 10004209
$ws.Range("I74").Value = 12501938
$ws.Range("J74").Value = 1678445.1
$ws.Range("K74").Value = 12501938
$ws.Range("L74").Value = 1678445.1
$ws.Range("M74").Value = -12501064
$ws.Range("N74").Value = -1680193.1
$ws.Range("H77").Value = 10004209
$ws.Range("I77").Value = 12501938
$ws.Range("J77").Value = 1678445.1
$ws.Range("K77").Value = 62509690
$ws.Range("L77").Value = 8392225.5
$ws.Range("M77").Value = -62505322
$ws.Range("N77").Value = -8400961.5
$ws.Range("H116").Value = 629.125
$ws.Range("I116").Value = 510
$ws.Range("J116").Value = 1224.75
$ws.Range("K116").Value = 510
$ws.Range("L116").Value = 1224.75
$ws.Range("M116").Value = 1784
$ws.Range("N116").Value = -5812.75
$ws.Range("H136").Value = 14710648
$ws.Range("I136").Value = 11115532
$ws.Range("J136").Value = 41674024
$ws.Range("K136").Value = 33346596
$ws.Range("L136").Value = 125022072
$ws.Range("M136").Value = -33344046
$ws.Range("N136").Value = -125027172

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 629.125
$ws.Range("I3").Value = 510
$ws.Range("J3").Value = 1224.75
$ws.Range("K3").Value = 510
$ws.Range("L3").Value = 1224.75
$ws.Range("M3").Value = -396
$ws.Range("N3").Value = -1452.75
$ws.Range("H22").Value = 2833.0715
$ws.Range("I22").Value = 2127.923
$ws.Range("J22").Value = 12000
$ws.Range("K22").Value = 2127.923
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = -1954.923
$ws.Range("N22").Value = -12346
$ws.Range("H94").Value = 880.25
$ws.Range("I94").Value = 742.1667
$ws.Range("J94").Value = 1018.3333
$ws.Range("K94").Value = 742.1667
$ws.Range("L94").Value = 1018.3333
$ws.Range("M94").Value = -291.1667
$ws.Range("N94").Value = -1920.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 42850
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 42850
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 42850
$ws.Range("N92").Value = -47842
$ws.Range("H107").Value = 1740.5518
$ws.Range("I107").Value = 750.875
$ws.Range("J107").Value = 2958.6155
$ws.Range("K107").Value = 750.875
$ws.Range("L107").Value = 2958.6155
$ws.Range("M107").Value = 1169.125
$ws.Range("N107").Value = -6798.6155
$ws.Range("H117").Value = 74036.336
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 74036.336
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 74036.336
$ws.Range("N117").Value = -83214.336
$ws.Range("H132").Value = 2309.1667
$ws.Range("I132").Value = 2319.1365
$ws.Range("J132").Value = 2199.5
$ws.Range("K132").Value = 6957.4095
$ws.Range("L132").Value = 6598.5
$ws.Range("M132").Value = -4427.4095
$ws.Range("N132").Value = -11658.5
$ws.Range("H134").Value = 2373.3076
$ws.Range("I134").Value = 1898.9688
$ws.Range("J134").Value = 4541.7144
$ws.Range("K134").Value = 5696.9064
$ws.Range("L134").Value = 13625.1432
$ws.Range("M134").Value = -3161.9064
$ws.Range("N134").Value = -18695.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2211.9092
$ws.Range("I5").Value = 2244.1
$ws.Range("J5").Value = 1890
$ws.Range("K5").Value = 6732.299999999999
$ws.Range("L5").Value = 5670
$ws.Range("M5").Value = -6620.299999999999
$ws.Range("N5").Value = -5894
$ws.Range("H22").Value = 1425
$ws.Range("I22").Value = 1425
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4275
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -4106
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 1425
$ws.Range("I27").Value = 1425
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 4275
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -4173
$ws.Range("N27").ClearContents()
$ws.Range("H68").Value = 1724.3334
$ws.Range("I68").Value = 1282
$ws.Range("J68").Value = 2166.6667
$ws.Range("K68").Value = 3846
$ws.Range("L68").Value = 6500.000100000001
$ws.Range("M68").Value = -3035
$ws.Range("N68").Value = -8122.000100000001
$ws.Range("H71").Value = 1724.3334
$ws.Range("I71").Value = 1282
$ws.Range("J71").Value = 2166.6667
$ws.Range("K71").Value = 11538
$ws.Range("L71").Value = 19500.0003
$ws.Range("M71").Value = -7482
$ws.Range("N71").Value = -27612.0003
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H135").Value = 2211.9092
$ws.Range("I135").Value = 2244.1
$ws.Range("J135").Value = 1890
$ws.Range("K135").Value = 20196.9
$ws.Range("L135").Value = 17010
$ws.Range("M135").Value = -17661.9
$ws.Range("N135").Value = -22080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1856.6666
$ws.Range("I97").Value = 1918.4286
$ws.Range("J97").Value = 1733.1428
$ws.Range("K97").Value = 1918.4286
$ws.Range("L97").Value = 1733.1428
$ws.Range("M97").Value = -1422.4286
$ws.Range("N97").Value = -2725.1428
$ws.Range("H132").Value = 16954004
$ws.Range("I132").Value = 23812092
$ws.Range("J132").Value = 10495.117
$ws.Range("K132").Value = 71436276
$ws.Range("L132").Value = 31485.351
$ws.Range("M132").Value = -71433746
$ws.Range("N132").Value = -36545.351

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 4000
$ws.Range("N4").Value = -4226
$ws.Range("H22").Value = 1679.8
$ws.Range("I22").Value = 1716
$ws.Range("J22").Value = 1625.5
$ws.Range("K22").Value = 1716
$ws.Range("L22").Value = 1625.5
$ws.Range("M22").Value = -1421
$ws.Range("N22").Value = -2215.5
$ws.Range("H27").Value = 1679.8
$ws.Range("I27").Value = 1716
$ws.Range("J27").Value = 1625.5
$ws.Range("K27").Value = 1716
$ws.Range("L27").Value = 1625.5
$ws.Range("M27").Value = -1609
$ws.Range("N27").Value = -1839.5
$ws.Range("H28").Value = 4000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 4000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 4000
$ws.Range("N28").Value = -4464
$ws.Range("H37").Value = 4000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 4000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 4000
$ws.Range("N37").Value = -4214
$ws.Range("H40").Value = 3268.4827
$ws.Range("I40").Value = 2128.7
$ws.Range("J40").Value = 3868.3684
$ws.Range("K40").Value = 2128.7
$ws.Range("L40").Value = 3868.3684
$ws.Range("M40").Value = -1992.7
$ws.Range("N40").Value = -4140.368399999999
$ws.Range("H46").Value = 2466.6428
$ws.Range("I46").Value = 2304.25
$ws.Range("J46").Value = 3441
$ws.Range("K46").Value = 2304.25
$ws.Range("L46").Value = 3441
$ws.Range("M46").Value = -2116.25
$ws.Range("N46").Value = -3817
$ws.Range("H93").Value = 100001270
$ws.Range("I93").Value = 250001300
$ws.Range("J93").Value = 1252.5
$ws.Range("K93").Value = 250001300
$ws.Range("L93").Value = 1252.5
$ws.Range("M93").Value = -250000052
$ws.Range("N93").Value = -3748.5
$ws.Range("H132").Value = 537453.9
$ws.Range("I132").Value = 14186.714
$ws.Range("J132").Value = 2002602
$ws.Range("K132").Value = 42560.142
$ws.Range("L132").Value = 6007806
$ws.Range("M132").Value = -40030.142
$ws.Range("N132").Value = -6012866

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 25231.846
$ws.Range("I15").Value = 17666.666
$ws.Range("J15").Value = 27501.4
$ws.Range("K15").Value = 17666.666
$ws.Range("L15").Value = 27501.4
$ws.Range("M15").Value = -17378.666
$ws.Range("N15").Value = -28077.4
$ws.Range("H21").Value = 10499.5
$ws.Range("I21").Value = 11000
$ws.Range("J21").Value = 9999
$ws.Range("K21").Value = 11000
$ws.Range("L21").Value = 9999
$ws.Range("M21").Value = -10765
$ws.Range("N21").Value = -10469
$ws.Range("H29").Value = 67500
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 67500
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 67500
$ws.Range("N29").Value = -68080
$ws.Range("H35").Value = 10499.5
$ws.Range("I35").Value = 11000
$ws.Range("J35").Value = 9999
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 9999
$ws.Range("M35").Value = -10710
$ws.Range("N35").Value = -10579
